# Boolean_Analysis_Output.xlsx edit
#
# For every worksheet:
#   - Drop the leading 0/1 "index" column (old column A) so the Yes/No
#     label column becomes column A and the Total column becomes column B.
#   - Add a new "Percentage" column C with a (count/total)*100 formula,
#     formatted to one decimal place, matching the bold/bordered header
#     style used by the other header cells.
#   - FM_Cooking_Demonstrations / FM_Kids_Activities also gained a second
#     data row (the "Yes"/0 row) that didn't previously appear.
#
# Column width + selection tweaks are reproduced on a best-effort basis to
# mirror the recorded diff; they are cosmetic only.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Farmers_Market_EBT  (total = 131 + 6 = 137)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

$ws.Columns("A").Delete()

$ws.Range("C1").Value = "Percentage"
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Font.Bold = $true

$ws.Range("C2").Formula = "=(B2/137)*100"
$ws.Range("C3").Formula = "=(B3/137)*100"
$ws.Range("C2:C3").NumberFormat = "0.0"

$ws.Range("C4").Select()

# ---------------------------------------------------------------------
# Sheet 2: FM_Year_Round  (total = 101 + 37 = 138)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Columns("A").Delete()

$ws.Range("C1").Value = "Percentage"
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Font.Bold = $true

$ws.Range("C2").Formula = "=(B2/138)*100"
$ws.Range("C3").Formula = "=(B3/138)*100"
$ws.Range("C2:C3").NumberFormat = "0.0"

$ws.Columns("A").ColumnWidth = 14.166666666666666
$ws.Columns("C").ColumnWidth = 9.666666666666666

$ws.Range("C3").Select()

# ---------------------------------------------------------------------
# Sheet 3: FM_Cooking_Demonstrations  (total = 138 + 0 = 138)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Columns("A").Delete()

# New second data row (previously missing).
$ws.Range("A3").Value = "Yes"
$ws.Range("B3").Value = 0

$ws.Range("C1").Value = "Percentage"
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Font.Bold = $true

$ws.Range("C2").Formula = "=(B2/138)*100"
$ws.Range("C3").Formula = "=(B3/138)*100"

$ws.Columns("A").ColumnWidth = 26.666666666666668
$ws.Columns("C").ColumnWidth = 9.666666666666666

$ws.Range("C3").Select()

# ---------------------------------------------------------------------
# Sheet 4: FM_Kids_Activities  (total = 138 + 0 = 138)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

$ws.Columns("A").Delete()

# New second data row (previously missing).
$ws.Range("A3").Value = "Yes"
$ws.Range("B3").Value = 0

$ws.Range("C1").Value = "Percentage"
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Font.Bold = $true

$ws.Range("C2").Formula = "=(B2/138)*100"
$ws.Range("C3").Formula = "=(B3/138)*100"

$ws.Columns("A").ColumnWidth = 20.666666666666668
$ws.Columns("B").ColumnWidth = 12.333333333333332
$ws.Columns("C").ColumnWidth = 9.333333333333332

$ws.Range("D4").Select()

# Re-select the first sheet/cell so the workbook opens where it did before.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

Write-Host "edit complete"
